$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Normalize capitalization of "de"/"la"/"las" in specific municipality names
$ws.Range("B3").Value = "San Cristóbal De Las Casas"
$ws.Range("B5").Value = "Jacala De Ledezma"
$ws.Range("B19").Value = "Martínez De La Torre"

# Remove trailing metadata/footer rows (24-28)
$ws.Range("A24:A28").EntireRow.Delete()
